# Re-saving the document with a newer Word build stamps the root element
# of several main parts with the w16sdtdh ("SDT data hash") markup-
# compatibility namespace, both as an xmlns declaration and as a token in
# mc:Ignorable. No actual body/content changes accompany this - it is a
# pure compatibility/namespace bump, so we round-trip the flat OOXML
# (WordOpenXML) and patch just the root tags of the affected parts.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$nsOld = 'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se='
$nsNew = 'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se='

$ignOldWp14 = 'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"'
$ignNewWp14 = 'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"'

$ignOldPlain = 'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex"'
$ignNewPlain = 'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh"'

$parts = @(
    "/word/document.xml",
    "/word/endnotes.xml",
    "/word/footer1.xml",
    "/word/footnotes.xml",
    "/word/header1.xml",
    "/word/numbering.xml",
    "/word/styles.xml"
)

foreach ($p in $parts) {
    $marker = 'pkg:name="' + $p + '"'
    $startIdx = $xml.IndexOf($marker)
    $dataOpen = $xml.IndexOf("<pkg:xmlData>", $startIdx) + 13
    $tagEnd = $xml.IndexOf(">", $dataOpen) + 1

    $before = $xml.Substring(0, $dataOpen)
    $rootTag = $xml.Substring($dataOpen, $tagEnd - $dataOpen)
    $after = $xml.Substring($tagEnd)

    $rootTag = $rootTag.Replace($nsOld, $nsNew)
    $rootTag = $rootTag.Replace($ignOldWp14, $ignNewWp14)
    $rootTag = $rootTag.Replace($ignOldPlain, $ignNewPlain)

    $xml = $before + $rootTag + $after
}

$d.Content.WordOpenXML = $xml
